$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 21:44"

# Update Estados Unidos (row 4)
$ws.Range("B4").Value = 7813715
$ws.Range("C4").Value = 37491
$ws.Range("D4").Value = 5008955
$ws.Range("E4").Value = 2587357
$ws.Range("G4").Value = 622
$ws.Range("H4").Value = 217403

# Update Alemania (row 26)
$ws.Range("B26").Value = 314947
$ws.Range("C26").Value = 3834
$ws.Range("E26").Value = 37581
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 9666

# Update Canada (row 29)
$ws.Range("B29").Value = 175067
$ws.Range("C29").Value = 1944
$ws.Range("D29").Value = 147218
$ws.Range("E29").Value = 18295

# Update Costa Rica (row 51)
$ws.Range("B51").Value = 84828
$ws.Range("C51").Value = 1331
$ws.Range("D51").Value = 51782
$ws.Range("E51").Value = 32006
$ws.Range("G51").Value = 16
$ws.Range("H51").Value = 1040

# Update Suiza (row 61)
$ws.Range("D61").Value = 48400
$ws.Range("E61").Value = 8394

# Update Tunez (row 82)
$ws.Range("B82").Value = 26899
$ws.Range("C82").Value = 2357
$ws.Range("E82").Value = 21458
$ws.Range("G82").Value = 45
$ws.Range("H82").Value = 409

# Update Maldivas (row 105)
$ws.Range("B105").Value = 10742
$ws.Range("C105").Value = 86
$ws.Range("D105").Value = 9589
$ws.Range("E105").Value = 1119

# Rows 123/124: Suazilandia moves above Lituania (new data for Suazilandia
# overtakes Lituania's unchanged totals), so swap country names and values
$ws.Range("A123").Value = "Suazilandia"
$ws.Range("B123").Value = 5632
$ws.Range("C123").Value = 15
$ws.Range("D123").Value = 5231
$ws.Range("E123").Value = 288
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 113

$ws.Range("A124").Value = "Lituania"
$ws.Range("B124").Value = 5625
$ws.Range("C124").Value = 142
$ws.Range("D124").Value = 2660
$ws.Range("E124").Value = 2863
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 102

# Update Sri Lanka (row 136)
$ws.Range("B136").Value = 4488
$ws.Range("C136").Value = 29
$ws.Range("E136").Value = 1197

# Update Aruba (row 138)
$ws.Range("B138").Value = 4150
$ws.Range("C138").Value = 17
$ws.Range("D138").Value = 3718
$ws.Range("E138").Value = 401

# Update Republica del Chad (row 166)
$ws.Range("B166").Value = 1262
$ws.Range("C166").Value = 11
$ws.Range("D166").Value = 1098
$ws.Range("E166").Value = 75
